$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 17: "Start Up Grant" program (Gruendungszuschuss) ---
# Caliendo et. al. (2015) - follow-up program to bridging allowance and
# start up subsidy, replacing both in 2006.
$ws.Range("A17").Value = "startupGrant"
$ws.Range("B17").Value = "Start Up Grant"
$ws.Range("C17").Value = 2009
$ws.Range("D17").Value = "Active Labor Market Policy"
$ws.Range("E17").Value = 40.969
$ws.Range("F17").Value = 'The start up Grant is the follow up program to the  bridging allowance "Überbrückungsgeld" and start up subsidy "Existenzgründungszuschuss" which were replaced in 2006. This subsidy pays the individual unemloyment benefit for 6 months and an additional 300 euros for up to 15 months.'
$ws.Range("G17").Value = 2047.962

# Match the wrap-text style used by the rest of column F/G
$ws.Range("F17:G17").WrapText = $true

# --- Row heights (recomputed by Excel for the wrapped long-text columns) ---
$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 75
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 90
$ws.Rows.Item(17).RowHeight = 90

# --- Selection state matches the author's saved view ---
[void]$ws.Range("A17").Select()
